$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191751003265381
$ws.Range("B1").Value = 2.425400733947754
$ws.Range("D1").Value = 1.428634643554688
$ws.Range("E1").Value = 0.9286888837814331
